{"js": "// Locate the paragraph that still contains the VORNAME/NACHNAME placeholder\n// (the \"Begr\u00fc\u00dfung\" intro paragraph) via search, then rewrite its whole text\n// with the real values substituted in, collapsing the paragraph down to a\n// single run - matching what Word does when Find&Replace consolidates a\n// paragraph made up of several same-formatted runs.\nconst body = context.document.body;\nconst results = body.search(\"VORNAME NACHNAME\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Placeholder paragraph not found\");\n}\n\nconst paragraph = results.items[0].paragraphs.getFirst();\nparagraph.load(\"text\");\nawait context.sync();\n\nlet newText = paragraph.text;\nnewText = newText.replace(\"VORNAME NACHNAME\", \"Samuel D\u00f6rr\");\nnewText = newText.replace(\"STERBEDATUM\", \"19/Mar/2021\");\nnewText = newText.replace(\"STERBEORT\", \"Duisburg\");\n\nparagraph.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Locate the paragraph that still contains the VORNAME/NACHNAME placeholder\n# (the \"Begr\u00fc\u00dfung\" intro paragraph), then rewrite its whole text with the\n# real values substituted in. Setting Range.Text collapses the paragraph's\n# several same-formatted runs into a single run - matching what Word does\n# when Find&Replace consolidates a paragraph like this.\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*VORNAME NACHNAME*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Placeholder paragraph not found\"\n}\n\n$para = $d.Paragraphs.Item($targetIndex)\n$r = $para.Range\n# Exclude the trailing paragraph mark from the range we rewrite.\n$r.End = $r.End - 1\n\n$newText = $r.Text\n$newText = $newText.Replace(\"VORNAME NACHNAME\", \"Samuel D\u00f6rr\")\n$newText = $newText.Replace(\"STERBEDATUM\", \"19/Mar/2021\")\n$newText = $newText.Replace(\"STERBEORT\", \"Duisburg\")\n\n$r.Text = $newText\n"}
